$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vinst & Vcham")

# Insert a new row at row 6 (pushes PP-Systems EGM-5 .. GAIA2TECH ECOFlux down by one,
# and Excel auto-updates the cross-sheet $C$13 -> $C$14 references on the "auxfile" sheet
# plus the in-sheet formula 45+B10+B11 -> 45+B11+B12).
$ws.Rows(6).Insert()

# Append the two new Aeris rows at the bottom of the table first, so the shared-string
# table receives "Aeris MIRA Ultra CH4/C2H6" / "Aeris MIRA Ultra N2O/CO2" before
# "Gasmet GT5000" (matches the order new strings were appended upstream).
$ws.Range("A16").Value = "Aeris MIRA Ultra CH4/C2H6"
$ws.Range("B16").Value = 60
$ws.Range("A17").Value = "Aeris MIRA Ultra N2O/CO2"
$ws.Range("B17").Value = 60

# Fill in the newly inserted row 6 with the Gasmet GT5000 instrument.
$ws.Range("A6").Value = "Gasmet GT5000"
$ws.Range("B6").Value = 500
$ws.Range("B6").Font.Color = 255

# Column A needs to widen to fit the longest label now in the table
# ("Aeris MIRA Ultra CH4/C2H6"), matching Excel's best-fit behaviour.
$ws.Columns("A").ColumnWidth = 22.83

# Match the recorded selection/active cell after the edit.
$ws.Range("B6").Select()
